# Implement quantitative metrics highlighting: split the plain-text run
# that contains an impact metric (percentage, dollar amount, large
# number) into multiple runs, applying bold + color (#2C3E50) to just
# the metric substrings.
#
# We operate paragraph-by-paragraph (via $d.Paragraphs) rather than a
# flat document-wide Find, because several of the target bullets are
# *prefixes* of other, longer bullets elsewhere in the resume (e.g. the
# "87% ... 71%" achievement also exists, with an extra trailing clause,
# under the "Partner - Siege Analytics" role) - a whole-document text
# search for the short form would otherwise match inside the long one.

function Find-ParaIndex {
    param($doc, $needle)
    $cnt = $doc.Paragraphs.Count
    for ($i = 1; $i -le $cnt; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# Paragraph range *excluding* the trailing paragraph-mark character, so
# InsertAfter never spills into the following paragraph.
function Get-ParaTextRange {
    param($doc, $idx)
    $p = $doc.Paragraphs.Item($idx)
    $pr = $p.Range
    $endpos = $pr.End - 1
    $r = $doc.Range($pr.Start, $endpos)
    return $r
}

function Add-PlainRun {
    param($range, $text)
    $range.InsertAfter($text)
    $range.Collapse(0)
}

function Add-MetricRun {
    param($range, $text)
    $range.InsertAfter($text)
    $range.Font.Bold = 1
    $range.Font.Color = 5258796
    $range.Collapse(0)
}

$d = $word.ActiveDocument
$bullet = [string][char]8226
$plusMinus = [string][char]177

# --- Paragraph: "Discovered systematic race coding errors ... from 23% to 64%"
$idx1 = Find-ParaIndex $d "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed"
if ($idx1 -gt 0) {
    $r = Get-ParaTextRange $d $idx1
    $r.Text = ""
    $seg = $bullet + " Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from "
    Add-PlainRun $r $seg
    Add-MetricRun $r "23%"
    $seg = " to "
    Add-PlainRun $r $seg
    Add-MetricRun $r "64%"
}

# --- Paragraph: "Achieved 87% ... 71%, reducing polling error margins from +-4.2% to +-2.1%"
$idx2 = Find-ParaIndex $d "reducing polling error margins"
if ($idx2 -gt 0) {
    $r = Get-ParaTextRange $d $idx2
    $r.Text = ""
    $seg = $bullet + " Achieved "
    Add-PlainRun $r $seg
    Add-MetricRun $r "87%"
    $seg = " prediction accuracy for voter turnout vs. industry standard of "
    Add-PlainRun $r $seg
    Add-MetricRun $r "71%"
    $seg = ", reducing polling error margins from "
    Add-PlainRun $r $seg
    $seg = $plusMinus + "4.2%"
    Add-MetricRun $r $seg
    $seg = " to "
    Add-PlainRun $r $seg
    $seg = $plusMinus + "2.1%"
    Add-MetricRun $r $seg
}

# --- Paragraph: "Wrote RFP and analyzed bids from 1,200 vendors ..."
$idx3 = Find-ParaIndex $d "Wrote RFP and analyzed bids from"
if ($idx3 -gt 0) {
    $r = Get-ParaTextRange $d $idx3
    $r.Text = ""
    $seg = $bullet + " Wrote RFP and analyzed bids from "
    Add-PlainRun $r $seg
    Add-MetricRun $r "1,200"
    $seg = " vendors for research platform development"
    Add-PlainRun $r $seg
}

# --- Paragraph: "Created comprehensive meta-analysis framework ... $400M ... $1B+"
$idx4 = Find-ParaIndex $d "Created comprehensive meta-analysis framework"
if ($idx4 -gt 0) {
    $r = Get-ParaTextRange $d $idx4
    $r.Text = ""
    $seg = $bullet + " Created comprehensive meta-analysis framework handling millions of survey responses that became the "
    Add-PlainRun $r $seg
    Add-MetricRun $r "`$400M"
    $seg = " Polling Consortium Database at The Analyst Institute, now valued at "
    Add-PlainRun $r $seg
    Add-MetricRun $r "`$1B"
    $seg = "+"
    Add-PlainRun $r $seg
}

# --- Paragraph: "Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M"
$idx5 = Find-ParaIndex $d "Algorithm reduced mapping costs"
if ($idx5 -gt 0) {
    $r = Get-ParaTextRange $d $idx5
    $r.Text = ""
    $seg = $bullet + " Algorithm reduced mapping costs by "
    Add-PlainRun $r $seg
    Add-MetricRun $r "73.5%"
    $seg = ", saving campaigns and organizations "
    Add-PlainRun $r $seg
    Add-MetricRun $r "`$4.7M"
}

# --- Paragraph: "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%" (short form)
$idx6 = -1
$cnt = $d.Paragraphs.Count
for ($i = 1; $i -le $cnt; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.Contains("Achieved 87%") -and (-not $t.Contains("reducing polling"))) {
        $idx6 = $i
    }
}
if ($idx6 -gt 0) {
    $r = Get-ParaTextRange $d $idx6
    $r.Text = ""
    $seg = $bullet + " Achieved "
    Add-PlainRun $r $seg
    Add-MetricRun $r "87%"
    $seg = " prediction accuracy for voter turnout vs. industry standard of "
    Add-PlainRun $r $seg
    Add-MetricRun $r "71%"
}

Write-Output "metrics highlighting applied"
